$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 3 and row 5 for columns D, J, K, L, M, N, O, P, Q
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

foreach ($col in $cols) {
    $addr3 = "$col" + "3"
    $addr5 = "$col" + "5"
    $val3 = $ws.Range($addr3).Value2
    $val5 = $ws.Range($addr5).Value2
    $ws.Range($addr3).Value = $val5
    $ws.Range($addr5).Value = $val3
}
